$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 12 ("Legislature" label row); everything below shifts up by one row.
$ws.Rows("12:12").Delete()

# Row 20 (new numbering): change C20:F20 from -0.5 to -0.2 (G20/H20 already -0.2)
$ws.Range("C20:F20").Value = -0.2

# Row 21 (new numbering): C21 and E21 get right-aligned (empty) style, matching style used elsewhere
$ws.Range("C21").HorizontalAlignment = -4152
$ws.Range("E21").HorizontalAlignment = -4152

# Update the selected/active cell to reflect where the user ended up editing
[void]$ws.Range("D25").Select()
